# Run modes for test suites
# - renames "makeDeposite" -> "depositeTest"
# - inserts a new "test_suite" sheet (TCID/Runmode table) between the
#   deposit sheet and "creditCalcCard"
# - leaves "creditCalcCard" content untouched (it just moves to the 3rd tab)

$wb = $excel.ActiveWorkbook

# 1) Rename the first sheet.
$depositeTest = $wb.Worksheets.Item(1)
$depositeTest.Name = "depositeTest"

# 2) Insert a brand-new sheet right after it (i.e. before creditCalcCard)
#    and name it "test_suite".
$testSuite = $wb.Worksheets.Add($null, $depositeTest)
$testSuite.Name = "test_suite"

# 3) Populate the run-mode table.
$testSuite.Range("A1").Value = "TCID"
$testSuite.Range("B1").Value = "Runmode"
$testSuite.Range("A2").Value = "CreditCalcCard"
$testSuite.Range("B2").Value = "Y"
$testSuite.Range("A3").Value = "DepositeTest"
$testSuite.Range("B3").Value = "Y"
$testSuite.Range("A4").Value = "OrderCardTest"
$testSuite.Range("B4").Value = "N"

# Match the look of the other data sheets (Calibri, same style as the
# creditCalcCard / depositeTest header rows).
$testSuite.Range("A1:B4").Font.Name = "Calibri"
